$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '67.828.88'
Set-TextValue $ws.Range('E2') '  -1.00%  '

Set-TextValue $ws.Range('D3') '3.788.55'
Set-TextValue $ws.Range('E3') '  +0.74%  '

Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  -0.14%  '

Set-TextValue $ws.Range('D5') '594.50'
Set-TextValue $ws.Range('E5') '  -0.34%  '

Set-TextValue $ws.Range('D6') '166.49'
Set-TextValue $ws.Range('E6') '  -1.10%  '

Set-TextValue $ws.Range('D7') '3.787.72'
Set-TextValue $ws.Range('E7') '  +0.90%  '

Set-TextValue $ws.Range('E8') '  +0.08%  '

Set-TextValue $ws.Range('E9') '  -0.71%  '

Set-TextValue $ws.Range('E10') '  -1.16%  '

Set-TextValue $ws.Range('E11') '  -2.01%  '

Set-TextValue $ws.Range('D12') '0.449'
Set-TextValue $ws.Range('E12') '  -0.54%  '

Set-TextValue $ws.Range('D13') '0.0000256'
Set-TextValue $ws.Range('E13') '  -2.04%  '

Set-TextValue $ws.Range('D14') '36.04'
Set-TextValue $ws.Range('E14') '  -1.32%  '

Set-TextValue $ws.Range('D15') '4.423.67'
Set-TextValue $ws.Range('E15') '  +0.64%  '

Set-TextValue $ws.Range('D16') '3.805.66'
Set-TextValue $ws.Range('E16') '  +1.13%  '

Set-TextValue $ws.Range('D17') '67.758.69'
Set-TextValue $ws.Range('E17') '  -1.21%  '

Set-TextValue $ws.Range('D18') '18.21'
Set-TextValue $ws.Range('E18') '  +0.93%  '

Set-TextValue $ws.Range('E19') '  +0.37%  '

Set-TextValue $ws.Range('D20') '6.99'
Set-TextValue $ws.Range('E20') '  -0.73%  '

Set-TextValue $ws.Range('D21') '10.25'
Set-TextValue $ws.Range('E21') '  -5.68%  '

Set-TextValue $ws.Range('D22') '460.40'
Set-TextValue $ws.Range('E22') '  -1.74%  '

Set-TextValue $ws.Range('D23') '0.695'
Set-TextValue $ws.Range('E23') '  -1.19%  '

Set-TextValue $ws.Range('E24') '  +3.56%  '

Set-TextValue $ws.Range('D25') '83.64'
Set-TextValue $ws.Range('E25') '  -1.04%  '

Set-TextValue $ws.Range('E26') '  -3.37%  '

Set-TextValue $ws.Range('D27') '11.88'
Set-TextValue $ws.Range('E27') '  -2.01%  '

Set-TextValue $ws.Range('D28') '10.11'
Set-TextValue $ws.Range('E28') '  -0.67%  '

Set-TextValue $ws.Range('E29') '  +0.14%  '

Set-TextValue $ws.Range('E30') '  -0.85%  '

$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D31') '7.26'
Set-TextValue $ws.Range('E31') '  -1.76%  '

$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D32') '29.87'
Set-TextValue $ws.Range('E32') '  -0.82%  '

Set-TextValue $ws.Range('D33') '2.20'
Set-TextValue $ws.Range('E33') '  -0.72%  '

Set-TextValue $ws.Range('D34') '9.16'
Set-TextValue $ws.Range('E34') '  -1.69%  '

Set-TextValue $ws.Range('D35') '0.998'
Set-TextValue $ws.Range('E35') '  -0.08%  '

Set-TextValue $ws.Range('D36') '3.741.41'
Set-TextValue $ws.Range('E36') '  +0.62%  '

Set-TextValue $ws.Range('D37') '0.100'
Set-TextValue $ws.Range('E37') '  -1.51%  '

Set-TextValue $ws.Range('D38') '3.34'
Set-TextValue $ws.Range('E38') '  -3.68%  '

Set-TextValue $ws.Range('E39') '  -0.46%  '

Set-TextValue $ws.Range('D40') '0.998'
Set-TextValue $ws.Range('E40') '  -0.83%  '

Set-TextValue $ws.Range('D41') '5.75'
Set-TextValue $ws.Range('E41') '  -1.62%  '

Set-TextValue $ws.Range('E42') '  -0.27%  '

Set-TextValue $ws.Range('E43') '  +0.01%  '

Set-TextValue $ws.Range('D44') '44.07'
Set-TextValue $ws.Range('E44') '  +1.08%  '

Set-TextValue $ws.Range('D45') '0.298'
Set-TextValue $ws.Range('E45') '  -3.16%  '

Set-TextValue $ws.Range('D46') '47.06'
Set-TextValue $ws.Range('E46') '  +2.34%  '

Set-TextValue $ws.Range('D47') '8.38'
Set-TextValue $ws.Range('E47') '  -2.76%  '

Set-TextValue $ws.Range('D48') '147.81'
Set-TextValue $ws.Range('E48') '  +1.05%  '

Set-TextValue $ws.Range('D49') '392.30'
Set-TextValue $ws.Range('E49') '  -1.41%  '

Set-TextValue $ws.Range('D50') '1.83'
Set-TextValue $ws.Range('E50') '  -6.68%  '

Set-TextValue $ws.Range('D51') '2.756.62'
Set-TextValue $ws.Range('E51') '  +2.51%  '
